$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($Row, $A, $D, $Values) {
    $ws.Cells.Item($Row, 1).Value = $A
    $ws.Cells.Item($Row, 4).Value = $D
    for ($i = 0; $i -lt $Values.Length; $i++) {
        $ws.Cells.Item($Row, 5 + $i).Value = $Values[$i]
    }
}

# Row 2: ECs -> Spon2 -> Itgam -> Resolving-Mac
Set-RowValues 2 "ECs" "Resolving-Mac" @(3,1,1.150782,3.452345999999999,0.03823856951930295,0.03823856951930295,3,1,35.68243999999999,107.04732,1,1,41.06270966807998,369.5643870127199,0.03823856951930295,0.03823856951930295)

# Row 3: FAPs -> Spon2 -> Itgam -> Resolving-Mac
Set-RowValues 3 "FAPs" "Resolving-Mac" @(3,1,27.59461233333333,82.78383700000001,0.9169230158851821,0.916923015885182,3,1,35.68243999999999,107.04732,1,1,984.6430989074265,8861.787890166839,0.9169230158851821,0.916923015885182)

# Row 4: MuSCs -> Spon2 -> Itgam -> Resolving-Mac
Set-RowValues 4 "MuSCs" "Resolving-Mac" @(3,1,1.290098666666667,3.870296,0.04286783035543951,0.0428678303554395,3,1,35.68243999999999,107.04732,1,1,46.03386826741332,414.3048144067199,0.04286783035543951,0.0428678303554395)

# Row 5: Resolving-Mac -> Spon2 -> Itgam -> Resolving-Mac
Set-RowValues 5 "Resolving-Mac" "Resolving-Mac" @(1,0.3333333333333333,0.05930433333333333,0.177913,0.001970584240075516,0.001970584240075516,3,1,35.68243999999999,107.04732,1,1,2.116123315906666,19.04510984316,0.001970584240075516,0.001970584240075516)

# Remove now-obsolete rows 6-9
$ws.Range("A6:T9").Delete()
